# Rotate the sighting data among rows 15, 16 and 17:
#   row15 -> row17, row16 -> row15, row17 -> row16
# (Row numbers / formatting stay put; only the record contents move.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture current ("before") values for the three rows ------------------
$row15 = @{
    A  = $ws.Range("A15").Value2
    B  = $ws.Range("B15").Value2
    E  = $ws.Range("E15").Value2
    F  = $ws.Range("F15").Value2
    G  = $ws.Range("G15").Value2
    H  = $ws.Range("H15").Value2
    Q  = $ws.Range("Q15").Value2
    R  = $ws.Range("R15").Value2
    AC = $ws.Range("AC15").Value2
    AX = $ws.Range("AX15").Value2
}

$row16 = @{
    A  = $ws.Range("A16").Value2
    B  = $ws.Range("B16").Value2
    E  = $ws.Range("E16").Value2
    F  = $ws.Range("F16").Value2
    G  = $ws.Range("G16").Value2
    H  = $ws.Range("H16").Value2
    Q  = $ws.Range("Q16").Value2
    R  = $ws.Range("R16").Value2
    AC = $ws.Range("AC16").Value2
    AX = $ws.Range("AX16").Value2
}

$row17 = @{
    A  = $ws.Range("A17").Value2
    B  = $ws.Range("B17").Value2
    E  = $ws.Range("E17").Value2
    F  = $ws.Range("F17").Value2
    G  = $ws.Range("G17").Value2
    H  = $ws.Range("H17").Value2
    Q  = $ws.Range("Q17").Value2
    R  = $ws.Range("R17").Value2
    AC = $ws.Range("AC17").Value2
    AX = $ws.Range("AX17").Value2
}

# --- write rotated ("after") values -----------------------------------------
# new row15 = old row16
$ws.Range("A15").Value = $row16.A
$ws.Range("B15").Value = $row16.B
$ws.Range("E15").Value = $row16.E
$ws.Range("F15").Value = $row16.F
$ws.Range("G15").Value = $row16.G
$ws.Range("H15").Value = $row16.H
$ws.Range("Q15").Value = $row16.Q
$ws.Range("R15").Value = $row16.R
if ($row16.AC) { $ws.Range("AC15").Value = $row16.AC } else { $ws.Range("AC15").ClearContents() }
$ws.Range("AX15").Value = $row16.AX

# new row16 = old row17
$ws.Range("A16").Value = $row17.A
$ws.Range("B16").Value = $row17.B
$ws.Range("E16").Value = $row17.E
$ws.Range("F16").Value = $row17.F
$ws.Range("G16").Value = $row17.G
$ws.Range("H16").Value = $row17.H
$ws.Range("Q16").Value = $row17.Q
$ws.Range("R16").Value = $row17.R
if ($row17.AC) { $ws.Range("AC16").Value = $row17.AC } else { $ws.Range("AC16").ClearContents() }
$ws.Range("AX16").Value = $row17.AX

# new row17 = old row15
$ws.Range("A17").Value = $row15.A
$ws.Range("B17").Value = $row15.B
$ws.Range("E17").Value = $row15.E
$ws.Range("F17").Value = $row15.F
$ws.Range("G17").Value = $row15.G
$ws.Range("H17").Value = $row15.H
$ws.Range("Q17").Value = $row15.Q
$ws.Range("R17").Value = $row15.R
if ($row15.AC) { $ws.Range("AC17").Value = $row15.AC } else { $ws.Range("AC17").ClearContents() }
$ws.Range("AX17").Value = $row15.AX
